# ASCOM Video Developers Documentation - updates from Hristo, 5 Nov 2013
$d = $word.ActiveDocument
$d.TrackRevisions = $false

# ---------------------------------------------------------------------------
# 1) Date updates: "October 2013" -> "November 2013" (title line + changelog)
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("October", $true, $false, $false, $false, $false, $true, 1, $false, "November", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Replace the "TODO" placeholder paragraph with a brand new section:
#    "Building Video Drivers for Digital Video Cameras"
# ---------------------------------------------------------------------------
$find2 = $d.Content.Find
$found = $find2.Execute("TODO: Still lots of things to add " + [char]0x2026, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$todoRange = $find2.Parent
$todoPara = $todoRange.Paragraphs(1)
$todoPara.Style = "Heading 1"
$todoRange.Text = "Building Video Drivers for Digital Video Cameras" + [char]9

# The paragraph right after the TODO line is already a lone empty paragraph
# in the source document (<w:p/>) - leave it untouched and insert the new
# content paragraphs after it.
$p2 = $todoPara.Next()

# Paragraph: intro text about dual mode cameras
$p2.Range.InsertParagraphAfter() | Out-Null
$p3 = $p2.Next()
$p3.Style = "Normal"
$p3.Range.Text = "Some digital video cameras can operate as both a video camera at a free-running mode and as a CCD camera in triggered exposure mode. In order to use the camera in both modes two separate drivers will be required - a Video driver and a Camera driver."

# Paragraph: explanation of video driver mode, with italic SupportedActions / Action()
$p3.Range.InsertParagraphAfter() | Out-Null
$p4 = $p3.Next()
$p4.Style = "Normal"
$p4.Range.Text = "When a client connects to the camera using the Video driver, then the camera needs to be set in a video mode, which will usually be a free-running video mode. If longer exposures in video mode are simulated via  individual exposures  triggered in CCD mode then this must be hidden from the client e.g. the exposures should be triggered by the driver automatically rather than by the client via the SupportedActions and Action() interface members. In all cases the video stream should be available at all times while the client is connected."

$f4 = $p4.Range.Find
$f4.Execute("SupportedActions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$f4.Parent.Font.Italic = $true

$f4b = $p4.Range.Find
$f4b.Execute("Action()", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$f4b.Parent.Font.Italic = $true

# Paragraph: Camera / Video interface exclusivity
$p4.Range.InsertParagraphAfter() | Out-Null
$p5 = $p4.Next()
$p5.Style = "Normal"
$p5.Range.Text = "If a client attempts to connect to the device using the Camera interface while there is also a client connected using the Video interface, then the second connection attempt must throw and error. A connection to each of the two supported interfaces must be possible only after all clients connected to the other interface have disconnected from the device. "

# Paragraph: page break before the Acknowledgment section.
# InsertBreak() on a collapsed range inside an *empty* paragraph always
# spawns an extra trailing empty paragraph in this runtime, so work around
# it: seed the new paragraph with a placeholder character, append the break
# after it (which merges the break into this same paragraph and produces a
# predictable single extra empty paragraph after it), strip the placeholder
# back out, then delete that extra trailing paragraph again.
$p5.Range.InsertParagraphAfter() | Out-Null
$p6 = $p5.Next()
$p6.Style = "Normal"
$p6.Range.Text = "X"
$breakRng = $p6.Range
$breakRng.Collapse(0)
$breakRng.InsertBreak(7) | Out-Null

$fx = $p6.Range.Find
$fx.Execute("X", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fx.Parent.Text = ""

$p7 = $p6.Next()
$p7.Range.Delete() | Out-Null
